# Auto-generated: refresh market-price derived columns (H-N) for specific Leve rows
# across all job sheets, per the scheduled market-data runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value = 3561.25
$ws.Range("I13").Value = 550
$ws.Range("J13").Value = 3991.4285
$ws.Range("K13").Value = 550
$ws.Range("L13").Value = 3991.4285
$ws.Range("M13").Value = -381
$ws.Range("N13").Value = -4329.4285
# Row 15
$ws.Range("H15").Value = 1454.96
$ws.Range("I15").Value = 1454.96
$ws.Range("K15").Value = 4364.88
$ws.Range("M15").Value = -4195.88
# Row 40
$ws.Range("H40").Value = 3283.4243
$ws.Range("J40").Value = 4292.857
$ws.Range("L40").Value = 4292.857
$ws.Range("N40").Value = -4642.857
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
# Row 62
$ws.Range("H62").Value = 4985.25
$ws.Range("I62").Value = 3978.3333
$ws.Range("K62").Value = 3978.3333
$ws.Range("M62").Value = -3354.3333
# Row 65
$ws.Range("H65").Value = 4985.25
$ws.Range("I65").Value = 3978.3333
$ws.Range("K65").Value = 19891.6665
$ws.Range("M65").Value = -16771.6665
# Row 88
$ws.Range("H88").Value = 623083.25
$ws.Range("I88").Value = 1958.8
$ws.Range("J88").Value = 1011286
$ws.Range("K88").Value = 1958.8
$ws.Range("L88").Value = 1011286
$ws.Range("M88").Value = -1552.8
$ws.Range("N88").Value = -1012098
# Row 91
$ws.Range("H91").Value = 623083.25
$ws.Range("I91").Value = 1958.8
$ws.Range("J91").Value = 1011286
$ws.Range("K91").Value = 1958.8
$ws.Range("L91").Value = 1011286
$ws.Range("M91").Value = -554.8
$ws.Range("N91").Value = -1014094
# Row 111
$ws.Range("H111").Value = 2406.4443
$ws.Range("I111").Value = 1972.25
$ws.Range("K111").Value = 5916.75
$ws.Range("M111").Value = -2849.75
# Row 133
$ws.Range("H133").Value = 63296.332
$ws.Range("J133").Value = 63296.332
$ws.Range("L133").Value = 63296.332
$ws.Range("N133").Value = -73416.33199999999
# Row 137
$ws.Range("H137").Value = 8551.375
$ws.Range("I137").Value = 3742
$ws.Range("K137").Value = 11226
$ws.Range("M137").Value = -8676
# Row 138
$ws.Range("H138").Value = 3100.7896
$ws.Range("J138").Value = 3344.9375
$ws.Range("L138").Value = 10034.8125
$ws.Range("N138").Value = -20314.8125

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1957.6
$ws.Range("I2").Value = 1957.6
$ws.Range("K2").Value = 1957.6
$ws.Range("M2").Value = -1844.6
# Row 32
$ws.Range("H32").Value = 10004357
$ws.Range("I32").Value = 10874294
$ws.Range("K32").Value = 10874294
$ws.Range("M32").Value = -10874007
# Row 61
$ws.Range("H61").Value = 27837740
$ws.Range("I61").Value = 45455910
$ws.Range("K61").Value = 45455910
$ws.Range("M61").Value = -45455698
# Row 116
$ws.Range("H116").Value = 1957.6
$ws.Range("I116").Value = 1957.6
$ws.Range("K116").Value = 1957.6
$ws.Range("M116").Value = 336.4000000000001
# Row 122
$ws.Range("H122").Value = 1412.9286
$ws.Range("I122").Value = 1252.6364
$ws.Range("K122").Value = 3757.9092
$ws.Range("M122").Value = -1307.9092
# Row 132
$ws.Range("H132").Value = 15035.375
$ws.Range("I132").Value = 9206
$ws.Range("K132").Value = 27618
$ws.Range("M132").Value = -25088
# Row 136
$ws.Range("H136").Value = 27837740
$ws.Range("I136").Value = 45455910
$ws.Range("K136").Value = 136367730
$ws.Range("M136").Value = -136365180

$ws = $wb.Worksheets.Item("BSM")
# Row 2
$ws.Range("H2").Value = 66989.8
$ws.Range("J2").Value = 66989.8
$ws.Range("L2").Value = 66989.8
$ws.Range("N2").Value = -67215.8
# Row 3
$ws.Range("H3").Value = 1957.6
$ws.Range("I3").Value = 1957.6
$ws.Range("K3").Value = 1957.6
$ws.Range("M3").Value = -1843.6
# Row 64
$ws.Range("H64").Value = 1479.625
$ws.Range("I64").Value = 1520
$ws.Range("K64").Value = 1520
$ws.Range("M64").Value = -1295
# Row 67
$ws.Range("H67").Value = 1479.625
$ws.Range("I67").Value = 1520
$ws.Range("K67").Value = 1520
$ws.Range("M67").Value = -740
# Row 95
$ws.Range("H95").Value = 18712.572
$ws.Range("J95").Value = 18712.572
$ws.Range("L95").Value = 18712.572
$ws.Range("N95").Value = -24204.572
# Row 105
$ws.Range("H105").Value = 1633.6666
$ws.Range("I105").Value = 1460.4
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1460.4
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = 286.5999999999999
$ws.Range("N105").Value = -5994
# Row 134
$ws.Range("H134").Value = 42133.9
$ws.Range("I134").Value = 6433.2915
$ws.Range("K134").Value = 19299.8745
$ws.Range("M134").Value = -16764.8745

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 840667.7
$ws.Range("I31").Value = 16608.9
$ws.Range("J31").Value = 1298478.1
$ws.Range("K31").Value = 16608.9
$ws.Range("L31").Value = 1298478.1
$ws.Range("M31").Value = -16313.9
$ws.Range("N31").Value = -1299068.1
# Row 34
$ws.Range("H34").Value = 840667.7
$ws.Range("I34").Value = 16608.9
$ws.Range("J34").Value = 1298478.1
$ws.Range("K34").Value = 16608.9
$ws.Range("L34").Value = 1298478.1
$ws.Range("M34").Value = -16406.9
$ws.Range("N34").Value = -1298882.1
# Row 86
$ws.Range("H86").Value = 4165.6665
$ws.Range("I86").Value = 3750
$ws.Range("K86").Value = 3750
$ws.Range("M86").Value = -2627
# Row 89
$ws.Range("H89").Value = 4165.6665
$ws.Range("I89").Value = 3750
$ws.Range("K89").Value = 18750
$ws.Range("M89").Value = -13134

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 997.8889
$ws.Range("J5").Value = 989
$ws.Range("L5").Value = 2967
$ws.Range("N5").Value = -3191
# Row 14
$ws.Range("H14").Value = 268.33334
$ws.Range("I14").Value = 268.33334
$ws.Range("K14").Value = 805.0000200000001
$ws.Range("M14").Value = -632.0000200000001
# Row 17
$ws.Range("H17").Value = 400
$ws.Range("J17").Value = 466.66666
$ws.Range("L17").Value = 1399.99998
$ws.Range("N17").Value = -1737.99998
# Row 34
$ws.Range("H34").Value = 2365
$ws.Range("J34").Value = 4499.8335
$ws.Range("L34").Value = 13499.5005
$ws.Range("N34").Value = -13667.5005
# Row 39
$ws.Range("H39").Value = 229390.44
$ws.Range("J39").Value = 219742.72
$ws.Range("L39").Value = 659228.16
$ws.Range("N39").Value = -659816.16
# Row 55
$ws.Range("H55").Value = 11175
$ws.Range("J55").Value = 11733.333
$ws.Range("L55").Value = 35199.999
$ws.Range("N55").Value = -35553.999
# Row 132
$ws.Range("H132").Value = 1700.75
$ws.Range("I132").Value = 1820.6666
$ws.Range("J132").Value = 1500.8889
$ws.Range("K132").Value = 16385.9994
$ws.Range("L132").Value = 13508.0001
$ws.Range("M132").Value = -13855.9994
$ws.Range("N132").Value = -18568.0001
# Row 135
$ws.Range("H135").Value = 997.8889
$ws.Range("J135").Value = 989
$ws.Range("L135").Value = 8901
$ws.Range("N135").Value = -13971

$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 14575.25
$ws.Range("I20").Value = 12400
$ws.Range("J20").Value = 16750.5
$ws.Range("K20").Value = 12400
$ws.Range("L20").Value = 16750.5
$ws.Range("M20").Value = -12155
$ws.Range("N20").Value = -17240.5
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
# Row 50
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
# Row 52
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
# Row 122
$ws.Range("H122").Value = 2749.1
$ws.Range("I122").Value = 2586.5
$ws.Range("K122").Value = 7759.5
$ws.Range("M122").Value = -5309.5
# Row 126
$ws.Range("H126").Value = 1996.5
$ws.Range("I126").Value = 1996.1428
$ws.Range("K126").Value = 5988.428400000001
$ws.Range("M126").Value = -3518.428400000001
# Row 132
$ws.Range("H132").Value = 71431980
$ws.Range("I132").Value = 100002790
$ws.Range("J132").Value = 4949.75
$ws.Range("K132").Value = 300008370
$ws.Range("L132").Value = 14849.25
$ws.Range("M132").Value = -300005840
$ws.Range("N132").Value = -19909.25

$ws = $wb.Worksheets.Item("LTW")
# Row 43
$ws.Range("H43").Value = 8767495
$ws.Range("I43").Value = 20000000
$ws.Range("J43").Value = 5023326.5
$ws.Range("K43").Value = 20000000
$ws.Range("L43").Value = 5023326.5
$ws.Range("M43").Value = -19999807
$ws.Range("N43").Value = -5023712.5
# Row 132
$ws.Range("H132").Value = 745058.7
$ws.Range("I132").Value = 771221.1
$ws.Range("K132").Value = 2313663.3
$ws.Range("M132").Value = -2311133.3

$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 30495
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
# Row 100
$ws.Range("H100").Value = 2046.4186
$ws.Range("I100").Value = 1999.8975
$ws.Range("K100").Value = 3999.795
$ws.Range("M100").Value = -3458.795
# Row 132
$ws.Range("H132").Value = 49999.8
$ws.Range("I132").Value = 49999
$ws.Range("K132").Value = 149997
$ws.Range("M132").Value = -147467
